# society import data feature improvement - area, subarea, city, state fixes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Supplier Code values fixed: KSH -> ASD, MIT -> ZXC
$ws.Range("E2").Value = "ASD"
$ws.Range("E3").Value = "ZXC"

# Row 3 height adjusted slightly
$ws.Rows.Item(3).RowHeight = 23.05

# Reset the view back to the top-left and move the active selection to E4
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E4").Select()
